$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the format of the last existing header cell (H1) onto the two new
# header cells so they share the same style index as the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Set the new header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Add new data values in columns I and J for rows 2-4
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9
$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9
